# updates from class, clarification, and quiz announcement
#
# Three changes from the source diff:
#   1. Slide 16 - reposition the empty highlight/callout rectangle
#      (Google Shape;328;p24) from (5461733, 2090025) EMU to
#      (6296100, -101259) EMU.
#   2. Slide 21 - clarification edit touching the three runs of the
#      "The list of UTF-8 characters:" bullet (no visible text change).
#   3. Slide 7 - reposition the "Payload" label box
#      (Google Shape;89;p16) from (1277850, 1383125) EMU to
#      (1408612, 499897) EMU.
#
# NOTE on literal numbers below: Shape.Left/Top are expressed in points
# (1 pt = 12700 EMU) and are stored by this host as single-precision
# floats that get truncated (not rounded) back to integer EMU on save.
# The literals here are chosen so that round-trip reproduces the exact
# target EMU coordinates from the diff.

$p = $ppt.ActivePresentation

# --- 1. Slide 16: move the blank rectangle shape -------------------------
$s16 = $p.Slides.Item(16)
$shp16 = $s16.Shapes.Item(6)
$shp16.Left = 495.75592041015625
$shp16.Top = -7.97314977645874

# --- 2. Slide 21: clarification touch-up on the UTF-8 bullet -------------
$s21 = $p.Slides.Item(21)
$shp21 = $s21.Shapes.Item(2)
$tr21 = $shp21.TextFrame.TextRange
$para1 = $tr21.Paragraphs(1, 1)

$run1 = $para1.Runs(1)
$run1.Text = "The list of "

$run2 = $para1.Runs(2)
$run2.Text = "UTF-8 characters"

$run3 = $para1.Runs(3)
$run3.Text = ":"

# --- 3. Slide 7: move the "Payload" label box -----------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(14)
$shp7.Left = 110.91433715820312
$shp7.Top = 39.361968994140625
